{"js": "// Trabalho2-LabES.docx: \"Novo trabalho de Lab. ES\"\n//\n// The title run \"Trabalho 1\" becomes \"Trabalho 2\", and Word's automatic\n// \"_GoBack\" bookmark (which marks the location of the last edit) moves from\n// the end of the document to right after the edited \"Trabalho 2\" run.\n\nconst body = context.document.body;\n\n// Locate the exact \"Trabalho 1\" run via search so we don't hard-code offsets.\nconst results = body.search(\"Trabalho 1\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('\"Trabalho 1\" not found in document body.');\n}\n\nconst titleRange = results.items[0];\n// Collapsed range immediately after \"Trabalho 1\" (before the following run).\nconst afterTitle = titleRange.getRange(\"After\");\n\n// Remove the existing \"_GoBack\" bookmark (currently at the end of the\n// document) and insert it right after the title text instead. Doing this\n// BEFORE editing the text keeps the run split at that boundary intact,\n// instead of letting the text edit re-merge it with the following run.\ncontext.document.deleteBookmark(\"_GoBack\");\nafterTitle.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// Now replace \"Trabalho 1\" with \"Trabalho 2\" in place.\ntitleRange.insertText(\"Trabalho 2\", \"Replace\");\nawait context.sync();\n", "ps1": "# Trabalho2-LabES.docx: \"Novo trabalho de Lab. ES\"\n#\n# The title run \"Trabalho 1\" becomes \"Trabalho 2\", and Word's automatic\n# \"_GoBack\" bookmark (which marks the location of the last edit) moves from\n# the end of the document to right after the edited \"Trabalho 2\" run.\n\n$d = $word.ActiveDocument\n\n# Locate the exact \"Trabalho 1\" run via Find so we don't hard-code offsets.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"Trabalho 1\")\nif (-not $found) {\n    throw \"'Trabalho 1' not found in document.\"\n}\n$titleStart = $findRange.Start\n$titleEnd = $findRange.End\n\n# Insert (relocate) the \"_GoBack\" bookmark immediately after the title text,\n# BEFORE editing the text, so the run split at that point is preserved\n# instead of being re-merged with the following run when the text changes.\n# Adding a bookmark with a name that already exists moves it, which also\n# takes care of removing the old \"_GoBack\" bookmark at the end of the doc.\n$bookmarkRange = $d.Range($titleEnd, $titleEnd)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange) | Out-Null\n\n# Now replace \"Trabalho 1\" with \"Trabalho 2\" in place.\n$titleRange = $d.Range($titleStart, $titleEnd)\n$titleRange.Text = \"Trabalho 2\"\n"}
